# Adds the "Best Practices" agenda-slide captions (slide 2) that were
# still empty placeholders in the deck. Each of the 7 numbered group
# boxes gets its caption text typed into its TextBox, colored with the
# deck's accent-blue (#004080), matching the language of its content.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$blue = 8404992  # RGB(0x00, 0x40, 0x80) -> R + G*256 + B*65536

function Set-RunText {
    param($TextBox, $Texts, $Langs)

    # Build the paragraph by inserting the chunks in *reverse* order,
    # each time via InsertBefore on the still-empty paragraph mark, so
    # the existing endParaRPr / pPr formatting is kept and the newly
    # inserted chunk is always run #1 at the moment we stamp its
    # language/colour onto it (the engine's LanguageID setter only
    # ever lands on the shape's first run, no matter which sub-range
    # it's invoked on - so we always call it while there is only one
    # run preceding the paragraph mark).
    for ($i = $Texts.Count - 1; $i -ge 0; $i--) {
        $tr = $TextBox.TextFrame.TextRange
        $tr.InsertBefore($Texts[$i])

        $full = $TextBox.TextFrame.TextRange
        $run = $full.Characters(1, $Texts[$i].Length)
        $run.Font.Color.RGB = $blue
        $run.LanguageID = $Langs[$i]
    }
}

# 1 -> Group 4 -> TextBox 26
Set-RunText $s.Shapes.Item(2).GroupItems.Item(3) @("Best Practices ") @("en-US")

# 2 -> Group 27 -> TextBox 30
Set-RunText $s.Shapes.Item(3).GroupItems.Item(3) @("Начинаем с наименее зависимых") @("ru-RU")

# 3 -> Group 31 -> TextBox 34
Set-RunText $s.Shapes.Item(4).GroupItems.Item(3) @("Простые ", "тесты") @("ru-RU", "ru-RU")

# 4 -> Group 35 -> TextBox 38
Set-RunText $s.Shapes.Item(5).GroupItems.Item(3) @("Константы в проверках") @("ru-RU")

# 6 -> Group 39 -> TextBox 42
Set-RunText $s.Shapes.Item(6).GroupItems.Item(3) @("Измеряйте покрытие") @("ru-RU")

# 5 -> Group 43 -> TextBox 46
Set-RunText $s.Shapes.Item(7).GroupItems.Item(3) @("Тестирование private ", "методов") @("it-IT", "it-IT")

# 7 -> Group 47 -> TextBox 50
Set-RunText $s.Shapes.Item(8).GroupItems.Item(3) @("Полная ", "автоматизиция") @("ru-RU", "ru-RU")
